$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates scraped from the crypto price/volume refresh.
# Columns D (Price) and E (Volume 1h) are plain text in this sheet (e.g. "29.214.50")
# so we force a text NumberFormat before assigning, then clear the format override
# back off (Excel would otherwise auto-parse the numeric-looking strings into floats).
$updates = [ordered]@{
    'D2' = '29.237.31'
    'E2' = '  -0.47%  '
    'D3' = '1.829.33'
    'E3' = '  -0.66%  '
    'D4' = '1.003'
    'E4' = '  +0.35%  '
    'D5' = '234.16'
    'E5' = '  -2.07%  '
    'D6' = '0.5981'
    'E6' = '  -4.68%  '
    'E7' = '  +0.43%  '
    'D8' = '0.06957'
    'E8' = '  -5.91%  '
    'D9' = '0.2743'
    'E9' = '  -5.13%  '
    'D10' = '23.25'
    'E10' = '  -6.46%  '
    'E11' = '  -1.18%  '
    'D12' = '1.824.41'
    'E12' = '  -1.11%  '
    'D13' = '4.751'
    'E13' = '  -4.35%  '
    'D14' = '0.6253'
    'E14' = '  -7.11%  '
    'D15' = '0.000009700'
    'E15' = '  -5.04%  '
    'D16' = '78.40'
    'E16' = '  -4.15%  '
    'D17' = '28.927.56'
    'E17' = '  -1.51%  '
    'D18' = '5.720'
    'E18' = '  -8.82%  '
    'D19' = '221.71'
    'E19' = '  -5.35%  '
    'E20' = '  +0.48%  '
    'D21' = '11.51'
    'E21' = '  -6.56%  '
    'D22' = '6.871'
    'E22' = '  -5.89%  '
    'D23' = '1.005'
    'E23' = '  +0.42%  '
    'D24' = '155.84'
    'E24' = '  -1.12%  '
    'D25' = '7.937'
    'E25' = '  -6.56%  '
    'D26' = '0.1289'
    'E26' = '  -4.02%  '
    'D27' = '16.49'
    'E27' = '  -4.71%  '
    'D28' = '0.06642'
    'E28' = '  -7.91%  '
    'D29' = '1.449'
    'E29' = '  -2.44%  '
    'D30' = '1.440'
    'E30' = '  -2.34%  '
    'D31' = '3.834'
    'E31' = '  -4.90%  '
    'D32' = '3.744'
    'E32' = '  -7.33%  '
    'D33' = '1.091'
    'E33' = '  -4.89%  '
    'D34' = '1.712'
    'E34' = '  -5.77%  '
    'D35' = '0.6420'
    'E35' = '  -8.13%  '
    'D36' = '2.547'
    'E36' = '  -0.99%  '
    'D37' = '2.735'
    'E37' = '  -2.50%  '
    'D38' = '1.180.37'
    'E38' = '  -4.28%  '
    'B39' = 'VeChain'
    'C39' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D39' = '0.01731'
    'E39' = '  -5.38%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D40' = '6.505'
    'E40' = '  -4.00%  '
    'D41' = '0.8994'
    'E41' = '  -5.17%  '
    'E42' = '  +0.45%  '
    'D43' = '1.978.45'
    'E43' = '  -0.75%  '
    'D44' = '100.44'
    'E44' = '  -0.60%  '
    'D45' = '62.01'
    'E45' = '  -4.99%  '
    'E46' = '  -2.97%  '
    'B47' = 'EnergySwap'
    'C47' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D47' = '8.464'
    'E47' = '  -4.67%  '
    'B48' = 'Cronos'
    'C48' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D48' = '0.05523'
    'E48' = '  -2.54%  '
    'D49' = '0.4553'
    'E49' = '  -0.49%  '
    'D50' = '1.574'
    'E50' = '  -7.34%  '
    'D51' = '6.328'
    'E51' = '  -9.07%  '
}

foreach ($ref in $updates.Keys) {
    $col = $ref -replace '[0-9]+$', ''
    $cell = $ws.Range($ref)
    if ($col -eq "D" -or $col -eq "E") {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$ref]
        $cell.ClearFormats()
    } else {
        $cell.Value = $updates[$ref]
    }
}
